$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("I2").Value = 0.08983953209358128
    $ws.Range("J2").Value = 0.0898395320935813
    $ws.Range("M2").Value = 0.794582
    $ws.Range("N2").Value = 2.383746
    $ws.Range("O2").Value = 0.03449752952410986
    $ws.Range("P2").Value = 0.03449752952410985
    $ws.Range("Q2").Value = 0.009519886942000001
    $ws.Range("R2").Value = 0.085678982478
    $ws.Range("S2").Value = 0.003099241910830536
    $ws.Range("T2").Value = 0.003099241910830535
    $ws.Range("I3").Value = 0.08983953209358128
    $ws.Range("J3").Value = 0.0898395320935813
    $ws.Range("O3").Value = 0.8945489325574519
    $ws.Range("P3").Value = 0.8945489325574517
    $ws.Range("S3").Value = 0.08036585753577408
    $ws.Range("T3").Value = 0.08036585753577408
    $ws.Range("I4").Value = 0.08983953209358128
    $ws.Range("J4").Value = 0.0898395320935813
    $ws.Range("M4").Value = 0.2871986666666667
    $ws.Range("N4").Value = 0.8615959999999999
    $ws.Range("O4").Value = 0.01246900191876775
    $ws.Range("P4").Value = 0.01246900191876775
    $ws.Range("Q4").Value = 0.003440927225333333
    $ws.Range("R4").Value = 0.030968345028
    $ws.Range("S4").Value = 0.001120209298056062
    $ws.Range("T4").Value = 0.001120209298056062
    $ws.Range("I5").Value = 0.08983953209358128
    $ws.Range("J5").Value = 0.0898395320935813
    $ws.Range("M5").Value = 1.149534666666667
    $ws.Range("N5").Value = 3.448604
    $ws.Range("O5").Value = 0.04990813547540859
    $ws.Range("P5").Value = 0.04990813547540859
    $ws.Range("Q5").Value = 0.01377257484133333
    $ws.Range("R5").Value = 0.123953173572
    $ws.Range("S5").Value = 0.004483723538773773
    $ws.Range("T5").Value = 0.004483723538773773
    $ws.Range("I6").Value = 0.08983953209358128
    $ws.Range("J6").Value = 0.0898395320935813
    $ws.Range("M6").Value = 0.1975403333333333
    $ws.Range("N6").Value = 0.5926210000000001
    $ws.Range("O6").Value = 0.008576400524262026
    $ws.Range("P6").Value = 0.008576400524262026
    $ws.Range("Q6").Value = 0.002366730733666667
    $ws.Range("R6").Value = 0.021300576603
    $ws.Range("S6").Value = 0.0007704998101468456
    $ws.Range("T6").Value = 0.0007704998101468457
    $ws.Range("G7").Value = 0.121379
    $ws.Range("H7").Value = 0.364137
    $ws.Range("I7").Value = 0.9101604679064187
    $ws.Range("J7").Value = 0.9101604679064187
    $ws.Range("M7").Value = 0.794582
    $ws.Range("N7").Value = 2.383746
    $ws.Range("O7").Value = 0.03449752952410986
    $ws.Range("P7").Value = 0.03449752952410985
    $ws.Range("Q7").Value = 0.096445568578
    $ws.Range("R7").Value = 0.8680101172019999
    $ws.Range("S7").Value = 0.03139828761327932
    $ws.Range("T7").Value = 0.03139828761327931
    $ws.Range("G8").Value = 0.121379
    $ws.Range("H8").Value = 0.364137
    $ws.Range("I8").Value = 0.9101604679064187
    $ws.Range("J8").Value = 0.9101604679064187
    $ws.Range("O8").Value = 0.8945489325574519
    $ws.Range("P8").Value = 0.8945489325574517
    $ws.Range("Q8").Value = 2.500911851124
    $ws.Range("R8").Value = 22.508206660116
    $ws.Range("S8").Value = 0.8141830750216779
    $ws.Range("T8").Value = 0.8141830750216776
    $ws.Range("G9").Value = 0.121379
    $ws.Range("H9").Value = 0.364137
    $ws.Range("I9").Value = 0.9101604679064187
    $ws.Range("J9").Value = 0.9101604679064187
    $ws.Range("M9").Value = 0.2871986666666667
    $ws.Range("N9").Value = 0.8615959999999999
    $ws.Range("O9").Value = 0.01246900191876775
    $ws.Range("P9").Value = 0.01246900191876775
    $ws.Range("Q9").Value = 0.03485988696133333
    $ws.Range("R9").Value = 0.313738982652
    $ws.Range("S9").Value = 0.01134879262071169
    $ws.Range("T9").Value = 0.01134879262071169
    $ws.Range("G10").Value = 0.121379
    $ws.Range("H10").Value = 0.364137
    $ws.Range("I10").Value = 0.9101604679064187
    $ws.Range("J10").Value = 0.9101604679064187
    $ws.Range("M10").Value = 1.149534666666667
    $ws.Range("N10").Value = 3.448604
    $ws.Range("O10").Value = 0.04990813547540859
    $ws.Range("P10").Value = 0.04990813547540859
    $ws.Range("Q10").Value = 0.1395293683053333
    $ws.Range("R10").Value = 1.255764314748
    $ws.Range("S10").Value = 0.04542441193663482
    $ws.Range("T10").Value = 0.04542441193663482
    $ws.Range("G11").Value = 0.121379
    $ws.Range("H11").Value = 0.364137
    $ws.Range("I11").Value = 0.9101604679064187
    $ws.Range("J11").Value = 0.9101604679064187
    $ws.Range("M11").Value = 0.1975403333333333
    $ws.Range("N11").Value = 0.5926210000000001
    $ws.Range("O11").Value = 0.008576400524262026
    $ws.Range("P11").Value = 0.008576400524262026
    $ws.Range("Q11").Value = 0.02397724811966667
    $ws.Range("R11").Value = 0.215795233077
    $ws.Range("S11").Value = 0.00780590071411518
    $ws.Range("T11").Value = 0.00780590071411518
